# Auto-generated from the cryptos.xlsx OOXML diff.
# Updates the Price (D) and Volume 1h (E) columns, plus the
# Avalanche/Chainlink row swap (B22:E23), to match the new snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.675.57"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.633.21"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.45"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("D12").Value = "1.861.65"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "1.675.36"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "26.672.36"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.53"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.22"
$ws.Range("E19").Value = "  +6.71%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.36"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  +5.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.65"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +3.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "1.226.52"
$ws.Range("E35").Value = "  +5.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +5.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.805"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "1.768.91"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.02"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.28"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("E50").Value = "  +4.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("E51").Value = "  -0.26%  "
